$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 98: politeness_score (column B) becomes a real number instead of text.
$ws.Cells.Item(98, 2).Value = 3

# Row 99: new annotation row appended below row 98.
$ws.Cells.Item(99, 1).Value = "Ruilin"
$ws.Cells.Item(99, 2).Value = "'3"
$ws.Cells.Item(99, 2).Style = "Normal"
$ws.Cells.Item(99, 3).Value = "无"
$ws.Cells.Item(99, 4).Value = "FBK"
$ws.Cells.Item(99, 5).Value = "EXP"
$ws.Cells.Item(99, 6).Value = "840f898f-6d0b-4603-abaa-7e0871215f61"
$ws.Cells.Item(99, 7).Value = "HyEi7bWR-_annotated.xlsx"
$ws.Cells.Item(99, 8).Value = "We have carried out additional experiments to examine run time and the following results will be included in the revision."
